$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Q2").Value = 1.95
$ws.Range("R2").Value = 1.9
$ws.Range("G3").Value = 1.65
$ws.Range("I3").Value = 6.25
$ws.Range("J3").Value = 2.3
$ws.Range("X3").Value = 6.5
$ws.Range("AA3").Value = 15
$ws.Range("AC3").Value = 7
$ws.Range("AD3").Value = 7
$ws.Range("AH3").Value = 29
$ws.Range("AI3").Value = 21
$ws.Range("AN3").Value = 3.4
$ws.Range("AV3").Value = 81
$ws.Range("G5").Value = 1.65
$ws.Range("H5").Value = 3.6
$ws.Range("J5").Value = 2.3
$ws.Range("O5").Value = 1.33
$ws.Range("P5").Value = 3.25
$ws.Range("Q5").Value = 2.05
$ws.Range("R5").Value = 1.75
$ws.Range("U5").Value = 2
$ws.Range("V5").Value = 1.73
$ws.Range("X5").Value = 7
$ws.Range("AC5").Value = 9
$ws.Range("AE5").Value = 19
$ws.Range("AF5").Value = 67
$ws.Range("AU5").Value = 9
$ws.Range("AV5").Value = 67
$ws.Range("Q7").Value = 3.4
$ws.Range("R7").Value = 1.33
$ws.Range("I12").Value = 4
$ws.Range("Y12").Value = 8.25
$ws.Range("AB12").Value = 22
$ws.Range("AD12").Value = 7.4
$ws.Range("AF12").Value = 55
$ws.Range("AG12").Value = 13.5
$ws.Range("AH12").Value = 24
$ws.Range("AJ12").Value = 60
$ws.Range("AX12").Value = 21
$ws.Range("G14").Value = 2.4
$ws.Range("I14").Value = 3
$ws.Range("K14").Value = 1.91
$ws.Range("AG14").Value = 7
$ws.Range("BA14").Value = 126
$ws.Range("O24").Value = 1.36
$ws.Range("P24").Value = 3
$ws.Range("Q24").Value = 2.15
$ws.Range("R24").Value = 1.67
$ws.Range("O25").Value = 1.29
$ws.Range("P25").Value = 3.5
$ws.Range("Q25").Value = 1.98
$ws.Range("R25").Value = 1.88
$ws.Range("G27").Value = 2.12
$ws.Range("I27").Value = 3.2
$ws.Range("J27").Value = 2.75
$ws.Range("K27").Value = 2.07
$ws.Range("L27").Value = 3.85
$ws.Range("M27").Value = 1.08
$ws.Range("N27").Value = 6.7
$ws.Range("O27").Value = 1.35
$ws.Range("P27").Value = 2.92
$ws.Range("Q27").Value = 2.05
$ws.Range("R27").Value = 1.72
$ws.Range("S27").Value = 1.44
$ws.Range("T27").Value = 2.62
$ws.Range("U27").Value = 1.82
$ws.Range("V27").Value = 1.88
$ws.Range("W27").Value = 7.1
$ws.Range("X27").Value = 10
$ws.Range("Y27").Value = 8.75
$ws.Range("Z27").Value = 20
$ws.Range("AA27").Value = 18
$ws.Range("AB27").Value = 30
$ws.Range("AC27").Value = 6.7
$ws.Range("AE27").Value = 15
$ws.Range("AF27").Value = 75
$ws.Range("AH27").Value = 16.5
$ws.Range("AI27").Value = 11.5
$ws.Range("AJ27").Value = 45
$ws.Range("AK27").Value = 30
$ws.Range("AL27").Value = 40
$ws.Range("AM27").Value = 600
$ws.Range("AN27").Value = 4
$ws.Range("AO27").Value = 11.25
$ws.Range("AT27").Value = 2.62
$ws.Range("AU27").Value = 7.2
$ws.Range("AV27").Value = 70
$ws.Range("AW27").Value = 5.2
$ws.Range("AX27").Value = 18.5
$ws.Range("AY27").Value = 26
$ws.Range("AZ27").Value = 100
$ws.Range("BA27").Value = 150
$ws.Range("BB27").Value = 350
$ws.Range("G28").Value = 2.4
$ws.Range("I28").Value = 2.77
$ws.Range("N28").Value = 6.7
$ws.Range("O28").Value = 1.34
$ws.Range("P28").Value = 3
$ws.Range("Q28").Value = 2.05
$ws.Range("R28").Value = 1.72
$ws.Range("T28").Value = 2.77
$ws.Range("U28").Value = 1.78
$ws.Range("V28").Value = 1.93
$ws.Range("W28").Value = 7.7
$ws.Range("X28").Value = 11.75
$ws.Range("Z28").Value = 26
$ws.Range("AA28").Value = 21
$ws.Range("AB28").Value = 30
$ws.Range("AC28").Value = 6.7
$ws.Range("AE28").Value = 14
$ws.Range("AF28").Value = 65
$ws.Range("AH28").Value = 14
$ws.Range("AI28").Value = 10.25
$ws.Range("AJ28").Value = 35
$ws.Range("AM28").Value = 500
$ws.Range("AN28").Value = 4.4
$ws.Range("AP28").Value = 20
$ws.Range("AS28").Value = 250
$ws.Range("AT28").Value = 2.77
$ws.Range("AU28").Value = 6.8
$ws.Range("AW28").Value = 4.8
$ws.Range("G29").Value = 2.57
$ws.Range("H29").Value = 3.45
$ws.Range("I29").Value = 2.37
$ws.Range("K29").Value = 2.18
$ws.Range("L29").Value = 2.95
$ws.Range("N29").Value = 7.9
$ws.Range("P29").Value = 3.55
$ws.Range("Q29").Value = 1.78
$ws.Range("R29").Value = 1.98
$ws.Range("S29").Value = 1.37
$ws.Range("T29").Value = 2.87
$ws.Range("W29").Value = 9.75
$ws.Range("X29").Value = 14
$ws.Range("Y29").Value = 9.75
$ws.Range("Z29").Value = 29
$ws.Range("AA29").Value = 20
$ws.Range("AC29").Value = 7.9
$ws.Range("AD29").Value = 6.8
$ws.Range("AF29").Value = 55
$ws.Range("AH29").Value = 12.5
$ws.Range("AJ29").Value = 25
$ws.Range("AK29").Value = 18.5
$ws.Range("AL29").Value = 26
$ws.Range("AN29").Value = 4.65
$ws.Range("AP29").Value = 21
$ws.Range("AT29").Value = 2.87
$ws.Range("AU29").Value = 7
$ws.Range("AV29").Value = 60
$ws.Range("AW29").Value = 4.45
$ws.Range("AX29").Value = 12.5
$ws.Range("AY29").Value = 19.5
$ws.Range("AZ29").Value = 50
$ws.Range("BA29").Value = 80
$ws.Range("BB29").Value = 250
$ws.Range("G30").Value = 1.87
$ws.Range("L30").Value = 3.65
$ws.Range("O30").Value = 1.15
$ws.Range("P30").Value = 4.7
$ws.Range("X30").Value = 12
$ws.Range("AL30").Value = 25
$ws.Range("AM30").Value = 200
$ws.Range("AY30").Value = 19
$ws.Range("BA30").Value = 80
$ws.Range("I33").Value = 2.15
$ws.Range("K33").Value = 2.22
$ws.Range("L33").Value = 2.72
$ws.Range("P33").Value = 3.55
$ws.Range("Q33").Value = 1.78
$ws.Range("R33").Value = 1.98
$ws.Range("S33").Value = 1.35
$ws.Range("T33").Value = 2.92
$ws.Range("X33").Value = 15.5
$ws.Range("AB33").Value = 30
$ws.Range("AK33").Value = 16.5
$ws.Range("AM33").Value = 400
$ws.Range("AN33").Value = 4.9
$ws.Range("AO33").Value = 15.5
$ws.Range("AP33").Value = 22
$ws.Range("AQ33").Value = 70
$ws.Range("AR33").Value = 100
$ws.Range("AT33").Value = 2.92
$ws.Range("AU33").Value = 7
$ws.Range("AV33").Value = 60
$ws.Range("AX33").Value = 11
$ws.Range("AY33").Value = 18.5
$ws.Range("BA33").Value = 70
$ws.Range("G35").Value = 1.65
$ws.Range("H35").Value = 3.4
$ws.Range("I35").Value = 5.3
$ws.Range("J35").Value = 2.27
$ws.Range("L35").Value = 5.2
$ws.Range("M35").Value = 1.02
$ws.Range("N35").Value = 9.75
$ws.Range("O35").Value = 1.27
$ws.Range("P35").Value = 3.1
$ws.Range("Q35").Value = 1.87
$ws.Range("R35").Value = 1.85
$ws.Range("T35").Value = 2.5
$ws.Range("U35").Value = 1.75
$ws.Range("V35").Value = 1.85
$ws.Range("W35").Value = 6.5
$ws.Range("X35").Value = 7.6
$ws.Range("Y35").Value = 7.9
$ws.Range("Z35").Value = 13
$ws.Range("AA35").Value = 13.5
$ws.Range("AC35").Value = 9.75
$ws.Range("AD35").Value = 6.7
$ws.Range("AE35").Value = 14.5
$ws.Range("AG35").Value = 14.5
$ws.Range("AH35").Value = 35
$ws.Range("AI35").Value = 16
$ws.Range("AJ35").Value = 110
$ws.Range("AK35").Value = 55
$ws.Range("AL35").Value = 50
$ws.Range("AN35").Value = 3.45
$ws.Range("AO35").Value = 8.25
$ws.Range("AQ35").Value = 29
$ws.Range("AU35").Value = 7.2
$ws.Range("AW35").Value = 6.7
$ws.Range("AX35").Value = 30
$ws.Range("AY35").Value = 32
$ws.Range("AZ35").Value = 200
$ws.Range("BA35").Value = 200
$ws.Range("BB35").Value = 400
$ws.Range("M37").Value = 1.02
$ws.Range("N37").Value = 7.1

Write-Host "Applied 227 cell updates"
